# Apply the LOM3205.xlsx update:
#  1. Change activation date 01/01/2019 -> 01/01/2023 (B8/C8)
#  2. Add English translation for "Objectives:" (B11/C11)
#  3. Add English translation for "Short syllabus:" (B13/C13)
#  4. Add English translation for "Syllabus:" (B15/C15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    # Route the literal text through a formula so Excel does not try to
    # auto-detect it as a date/number, then bake the formula down to a
    # plain value via copy / paste-values so the cell ends up as a normal
    # shared-string cell (not a formula cell).
    $escaped = $Text.Replace("""", """""")
    $Range.Formula = "=""" + $escaped + """"
    $Range.Copy() | Out-Null
    $Range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# --- 1. Activation date ------------------------------------------------
Set-TextValue $ws.Range("B8") "01/01/2023"
Set-TextValue $ws.Range("C8") "01/01/2023"

# --- 2. Objectives (English) -------------------------------------------
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

Set-TextValue $ws.Range("B11") "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."
Set-TextValue $ws.Range("C11") "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."

# --- 3. Short syllabus (English) ---------------------------------------
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

Set-TextValue $ws.Range("B13") "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"
Set-TextValue $ws.Range("C13") "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"

# --- 4. Syllabus (English) ----------------------------------------------
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

Set-TextValue $ws.Range("B15") "Electrostatics (electrostatic field; electric potential; work and energy in electrostatics). Special techniques for solving the Laplace’s equation (method of images; separation of variables). Electric field in matter (electric polarization; polarized object field; bound charges; electric displacement; dielectrics (linear). Magnetostatics (Lorentz's law; Biot-Savart's law; Ampere's law; vector magnetic potential). Magnetic field in matter (magnetization; field of a magnetized object; auxiliary field H). Electrodynamics (electromotive force; electromagnetic induction; Maxwell's equations; law of conservation of charge). Electromagnetic waves (propagation in vacuum and in matter; reflection and transmission), plane wave equation and boundary conditions (interfaces). Electric dipole radiation."
Set-TextValue $ws.Range("C15") "Electrostatics (electrostatic field; electric potential; work and energy in electrostatics). Special techniques for solving the Laplace’s equation (method of images; separation of variables). Electric field in matter (electric polarization; polarized object field; bound charges; electric displacement; dielectrics (linear). Magnetostatics (Lorentz's law; Biot-Savart's law; Ampere's law; vector magnetic potential). Magnetic field in matter (magnetization; field of a magnetized object; auxiliary field H). Electrodynamics (electromotive force; electromagnetic induction; Maxwell's equations; law of conservation of charge). Electromagnetic waves (propagation in vacuum and in matter; reflection and transmission), plane wave equation and boundary conditions (interfaces). Electric dipole radiation."

$ws.Range("A1").Select() | Out-Null

Write-Output "done"
